$d = $word.ActiveDocument

# Update the date in the title paragraph
$d.Content.Find.Execute("2023-10-10 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-11 Wednesday", 2) | Out-Null

# Update the division problems/answers in the table, addressed by cell
# to avoid collisions between old and new values across different cells.
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "33÷2=16, 1"
$t.Cell(1, 2).Range.Text = "94÷6=15, 4"
$t.Cell(1, 3).Range.Text = "61÷5=12, 1"
$t.Cell(1, 4).Range.Text = "42÷3=14, 0"
$t.Cell(1, 5).Range.Text = "65÷3=21, 2"
$t.Cell(5, 1).Range.Text = "53÷9=5, 8"
$t.Cell(5, 2).Range.Text = "28÷6=4, 4"
$t.Cell(5, 3).Range.Text = "90÷4=22, 2"
$t.Cell(5, 4).Range.Text = "46÷3=15, 1"
$t.Cell(5, 5).Range.Text = "48÷2=24, 0"
$t.Cell(9, 1).Range.Text = "83÷3=27, 2"
$t.Cell(9, 2).Range.Text = "66÷9=7, 3"
$t.Cell(9, 3).Range.Text = "93÷8=11, 5"
$t.Cell(9, 4).Range.Text = "61÷6=10, 1"
$t.Cell(9, 5).Range.Text = "39÷7=5, 4"
$t.Cell(13, 1).Range.Text = "68÷5=13, 3"
$t.Cell(13, 2).Range.Text = "33÷8=4, 1"
$t.Cell(13, 3).Range.Text = "81÷4=20, 1"
$t.Cell(13, 4).Range.Text = "36÷4=9, 0"
$t.Cell(13, 5).Range.Text = "82÷4=20, 2"
$t.Cell(17, 1).Range.Text = "11÷8=1, 3"
$t.Cell(17, 2).Range.Text = "88÷8=11, 0"
$t.Cell(17, 3).Range.Text = "71÷3=23, 2"
$t.Cell(17, 4).Range.Text = "58÷8=7, 2"
$t.Cell(17, 5).Range.Text = "78÷8=9, 6"
